# Normalize the "Recorded By" (column G) values on the "Session Analysis
# Results" sheet: a handful of specific, exact string values get their
# comma-separated tokens reordered (effectively a left-rotation moving the
# first token to the end). Only cells whose full text matches one of the
# known source strings exactly are touched; everything else is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Exact-match lookup table: old full cell text -> new full cell text.
$replacements = @{
    "System, system, backup@backdoor.com" = "system, backup@backdoor.com, System"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "System, backup@backdoor.com"         = "backup@backdoor.com, System"
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value2
    if ($null -ne $val -and $replacements.ContainsKey($val)) {
        $cell.Value2 = $replacements[$val]
    }
}
